# Apply corrected fixed-recourse data and new MP iteration rows
# (per commit: "changed MP time limit and corrected error in fixed recourse data")
$wb = $excel.ActiveWorkbook

# ---- Sheet1: per-instance summary table ----
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Cells.Item(2, 2).Value = -274.9655817221359
$ws1.Cells.Item(2, 3).Value = 9.011785872
$ws1.Cells.Item(2, 4).Value = 3
$ws1.Cells.Item(2, 6).Value = 2
$ws1.Cells.Item(2, 7).Value = 1204
$ws1.Cells.Item(2, 8).Value = 1210
$ws1.Cells.Item(2, 9).Value = 100
$ws1.Cells.Item(3, 2).Value = -273.98162258884673
$ws1.Cells.Item(3, 3).Value = 5.336540295
$ws1.Cells.Item(3, 4).Value = 6
$ws1.Cells.Item(3, 6).Value = 5
$ws1.Cells.Item(3, 7).Value = 2860
$ws1.Cells.Item(3, 8).Value = 3025
$ws1.Cells.Item(3, 9).Value = 250
$ws1.Cells.Item(4, 2).Value = -274.08960459636427
$ws1.Cells.Item(4, 3).Value = 6.191085257
$ws1.Cells.Item(4, 4).Value = 7
$ws1.Cells.Item(4, 6).Value = 6
$ws1.Cells.Item(4, 7).Value = 3412
$ws1.Cells.Item(4, 8).Value = 3630
$ws1.Cells.Item(4, 9).Value = 300
$ws1.Cells.Item(5, 2).Value = -276.8685515416252
$ws1.Cells.Item(5, 3).Value = 4.272633059
$ws1.Cells.Item(5, 4).Value = 5
$ws1.Cells.Item(5, 6).Value = 4
$ws1.Cells.Item(5, 7).Value = 2308
$ws1.Cells.Item(5, 8).Value = 2420
$ws1.Cells.Item(5, 9).Value = 200
$ws1.Cells.Item(6, 2).Value = -272.1653938311721
$ws1.Cells.Item(6, 3).Value = 12.95195401
$ws1.Cells.Item(6, 4).Value = 8
$ws1.Cells.Item(6, 6).Value = 7
$ws1.Cells.Item(6, 7).Value = 3964
$ws1.Cells.Item(6, 8).Value = 4235
$ws1.Cells.Item(6, 9).Value = 350
$ws1.Cells.Item(7, 2).Value = -268.97221187709164
$ws1.Cells.Item(7, 3).Value = 3.384116984
$ws1.Cells.Item(7, 4).Value = 4
$ws1.Cells.Item(7, 6).Value = 3
$ws1.Cells.Item(7, 7).Value = 1756
$ws1.Cells.Item(7, 8).Value = 1815
$ws1.Cells.Item(7, 9).Value = 150
$ws1.Cells.Item(8, 2).Value = -265.42815137347486
$ws1.Cells.Item(8, 3).Value = 2.785066582
$ws1.Cells.Item(8, 4).Value = 2
$ws1.Cells.Item(8, 6).Value = 1
$ws1.Cells.Item(8, 7).Value = 652
$ws1.Cells.Item(8, 8).Value = 605
$ws1.Cells.Item(8, 9).Value = 50
$ws1.Cells.Item(9, 2).Value = -274.20170847509786
$ws1.Cells.Item(9, 3).Value = 1.601040059
$ws1.Cells.Item(9, 4).Value = 3
$ws1.Cells.Item(9, 6).Value = 2
$ws1.Cells.Item(9, 7).Value = 1204
$ws1.Cells.Item(9, 8).Value = 1210
$ws1.Cells.Item(9, 9).Value = 100
$ws1.Cells.Item(10, 2).Value = -271.53604073578464
$ws1.Cells.Item(10, 3).Value = 10.289141947
$ws1.Cells.Item(10, 4).Value = 7
$ws1.Cells.Item(10, 6).Value = 6
$ws1.Cells.Item(10, 7).Value = 3412
$ws1.Cells.Item(10, 8).Value = 3630
$ws1.Cells.Item(10, 9).Value = 300
$ws1.Cells.Item(11, 2).Value = -268.78676349663795
$ws1.Cells.Item(11, 3).Value = 1.26981565
$ws1.Cells.Item(11, 4).Value = 3
$ws1.Cells.Item(11, 6).Value = 2
$ws1.Cells.Item(11, 7).Value = 1204
$ws1.Cells.Item(11, 8).Value = 1210
$ws1.Cells.Item(11, 9).Value = 100

# ---- Sheet "1": master-problem iteration log ----
$ws = $wb.Worksheets.Item("1")
$ws.Cells.Item(2, 4).Value = 0.8459882274916992
$ws.Cells.Item(2, 5).Value = 70.77445
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = -275.37819416649336
$ws.Cells.Item(3, 3).Value = 0.08776712157840491
$ws.Cells.Item(3, 4).Value = 0.1868246894711914
$ws.Cells.Item(3, 5).Value = 0.93488
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = -274.9655817221359
$ws.Cells.Item(4, 3).Value = 0.08770811021198437
$ws.Cells.Item(4, 4).Value = 0.42188376053234866
$ws.Cells.Item(4, 5).Value = 0.0

# ---- Sheet "2": master-problem iteration log ----
$ws = $wb.Worksheets.Item("2")
$ws.Cells.Item(2, 4).Value = 0.043430891670776364
$ws.Cells.Item(2, 5).Value = 68.31078
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = -282.5400242680931
$ws.Cells.Item(3, 3).Value = 0.09620453204212696
$ws.Cells.Item(3, 4).Value = 0.5297453552442627
$ws.Cells.Item(3, 5).Value = 2.03306
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = -279.9858049039651
$ws.Cells.Item(4, 3).Value = 0.026970352815420248
$ws.Cells.Item(4, 4).Value = 0.77098059350708
$ws.Cells.Item(4, 5).Value = 2.06672
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = -276.7415271190988
$ws.Cells.Item(5, 3).Value = 0.09205844327296142
$ws.Cells.Item(5, 4).Value = 0.8785558021623535
$ws.Cells.Item(5, 5).Value = 1.71613
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = -274.93974763381715
$ws.Cells.Item(6, 3).Value = 0.0800621075189235
$ws.Cells.Item(6, 4).Value = 1.0210462858820801
$ws.Cells.Item(6, 5).Value = 1.03583
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = -273.98162258884673
$ws.Cells.Item(7, 3).Value = 0.08373334482089792
$ws.Cells.Item(7, 4).Value = 1.6440672093690185
$ws.Cells.Item(7, 5).Value = 0.0

# ---- Sheet "3": master-problem iteration log ----
$ws = $wb.Worksheets.Item("3")
$ws.Cells.Item(2, 4).Value = 0.019674989152832032
$ws.Cells.Item(2, 5).Value = 72.03609
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = -284.93216461805935
$ws.Cells.Item(3, 3).Value = 0.06856170788369287
$ws.Cells.Item(3, 4).Value = 0.3585616605776367
$ws.Cells.Item(3, 5).Value = 2.6296
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = -281.02342668428673
$ws.Cells.Item(4, 3).Value = 0.014996801358430454
$ws.Cells.Item(4, 4).Value = 0.5283823002906494
$ws.Cells.Item(4, 5).Value = 1.78881
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = -275.8641653669497
$ws.Cells.Item(5, 3).Value = 0.0975185014249835
$ws.Cells.Item(5, 4).Value = 0.7497194855341797
$ws.Cells.Item(5, 5).Value = 1.34265
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = -274.76650879323324
$ws.Cells.Item(6, 3).Value = 0.0
$ws.Cells.Item(6, 4).Value = 1.1823969776972656
$ws.Cells.Item(6, 5).Value = 0.68833
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = -274.1144998692482
$ws.Cells.Item(7, 3).Value = 0.052653302771760724
$ws.Cells.Item(7, 4).Value = 1.0471213507736816
$ws.Cells.Item(7, 5).Value = 0.20688
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = -274.08960459636427
$ws.Cells.Item(8, 3).Value = 0.07708650067600377
$ws.Cells.Item(8, 4).Value = 1.7273631935544433
$ws.Cells.Item(8, 5).Value = 0.0

# ---- Sheet "4": master-problem iteration log ----
$ws = $wb.Worksheets.Item("4")
$ws.Cells.Item(2, 4).Value = 0.03315735772436523
$ws.Cells.Item(2, 5).Value = 71.29895
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = -278.300849315347
$ws.Cells.Item(3, 3).Value = 0.09910138944675462
$ws.Cells.Item(3, 4).Value = 0.5236881781154785
$ws.Cells.Item(3, 5).Value = 2.02028
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = -277.44746019600944
$ws.Cells.Item(4, 3).Value = 0.09324511586117006
$ws.Cells.Item(4, 4).Value = 0.8548955550218506
$ws.Cells.Item(4, 5).Value = 1.24061
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = -276.95668727849005
$ws.Cells.Item(5, 3).Value = 0.09548483959325342
$ws.Cells.Item(5, 4).Value = 1.029056386616333
$ws.Cells.Item(5, 5).Value = 0.59438
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = -276.8685515416252
$ws.Cells.Item(6, 3).Value = 0.0982522559478422
$ws.Cells.Item(6, 4).Value = 1.4925861974378662
$ws.Cells.Item(6, 5).Value = 0.0

# ---- Sheet "5": master-problem iteration log ----
$ws = $wb.Worksheets.Item("5")
$ws.Cells.Item(2, 4).Value = 0.05066696337133789
$ws.Cells.Item(2, 5).Value = 70.50756
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = -321.98521738903344
$ws.Cells.Item(3, 3).Value = 0.05901927250361406
$ws.Cells.Item(3, 4).Value = 0.11446679339587403
$ws.Cells.Item(3, 5).Value = 5.48806
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = -297.0793900691263
$ws.Cells.Item(4, 3).Value = 0.0
$ws.Cells.Item(4, 4).Value = 0.7611371963791503
$ws.Cells.Item(4, 5).Value = 3.9308
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = -282.40881917382546
$ws.Cells.Item(5, 3).Value = 0.0
$ws.Cells.Item(5, 4).Value = 0.985951380078003
$ws.Cells.Item(5, 5).Value = 2.64459
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = -277.1180357760685
$ws.Cells.Item(6, 3).Value = 0.08381775144543133
$ws.Cells.Item(6, 4).Value = 0.9481974394178467
$ws.Cells.Item(6, 5).Value = 1.93829
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = -273.1943884633127
$ws.Cells.Item(7, 3).Value = 0.07907803844069364
$ws.Cells.Item(7, 4).Value = 0.9899204675969239
$ws.Cells.Item(7, 5).Value = 0.91228
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = -272.5171413614702
$ws.Cells.Item(8, 3).Value = 0.08563083692252649
$ws.Cells.Item(8, 4).Value = 1.3552736605531006
$ws.Cells.Item(8, 5).Value = 0.43816
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = -272.1653938311721
$ws.Cells.Item(9, 3).Value = 0.07516289901403857
$ws.Cells.Item(9, 4).Value = 6.999897805318115
$ws.Cells.Item(9, 5).Value = 0.0

# ---- Sheet "6": master-problem iteration log ----
$ws = $wb.Worksheets.Item("6")
$ws.Cells.Item(2, 4).Value = 0.03487423657678223
$ws.Cells.Item(2, 5).Value = 74.8484
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = -272.4220398086923
$ws.Cells.Item(3, 3).Value = 0.002229114480537501
$ws.Cells.Item(3, 4).Value = 0.8938715389545898
$ws.Cells.Item(3, 5).Value = 1.64852
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = -268.9931632473488
$ws.Cells.Item(4, 3).Value = 0.09418428710681696
$ws.Cells.Item(4, 4).Value = 1.2094037405688476
$ws.Cells.Item(4, 5).Value = 0.24472
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = -268.97221187709164
$ws.Cells.Item(5, 3).Value = 0.09712584067328685
$ws.Cells.Item(5, 4).Value = 0.9518760947789306
$ws.Cells.Item(5, 5).Value = 0.0

# ---- Sheet "7": master-problem iteration log ----
$ws = $wb.Worksheets.Item("7")
$ws.Cells.Item(2, 4).Value = 0.045845211184448245
$ws.Cells.Item(2, 5).Value = 69.60139
$ws.Cells.Item(3, 2).Value = -265.42815137347486
$ws.Cells.Item(3, 3).Value = 0.002468932748980447
$ws.Cells.Item(3, 4).Value = 2.637062962984253

# ---- Sheet "8": master-problem iteration log ----
$ws = $wb.Worksheets.Item("8")
$ws.Cells.Item(2, 4).Value = 0.052811610353881835
$ws.Cells.Item(2, 5).Value = 71.06337
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = -274.9067546246747
$ws.Cells.Item(3, 3).Value = 0.009771497210873632
$ws.Cells.Item(3, 4).Value = 0.20051590111315917
$ws.Cells.Item(3, 5).Value = 1.14379
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = -274.20170847509786
$ws.Cells.Item(4, 3).Value = 0.09846319087907521
$ws.Cells.Item(4, 4).Value = 1.1591712983615723
$ws.Cells.Item(4, 5).Value = 0.0

# ---- Sheet "9": master-problem iteration log ----
$ws = $wb.Worksheets.Item("9")
$ws.Cells.Item(2, 4).Value = 0.03476508155932617
$ws.Cells.Item(2, 5).Value = 68.71954
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = -306.4514620630372
$ws.Cells.Item(3, 3).Value = 0.007815128475927568
$ws.Cells.Item(3, 4).Value = 0.20442813810375976
$ws.Cells.Item(3, 5).Value = 4.76064
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = -282.48717157721865
$ws.Cells.Item(4, 3).Value = 0.07148866119377287
$ws.Cells.Item(4, 4).Value = 1.8136089265018311
$ws.Cells.Item(4, 5).Value = 2.74184
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = -278.1525706316617
$ws.Cells.Item(5, 3).Value = 0.06081703560710542
$ws.Cells.Item(5, 4).Value = 0.7090817878721923
$ws.Cells.Item(5, 5).Value = 2.11724
$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = -271.91021281294024
$ws.Cells.Item(6, 3).Value = 0.056022989450732835
$ws.Cells.Item(6, 4).Value = 1.4961259629753418
$ws.Cells.Item(6, 5).Value = 1.21104
$ws.Cells.Item(7, 1).Value = 6
$ws.Cells.Item(7, 2).Value = -271.7016546746654
$ws.Cells.Item(7, 3).Value = 0.09421008177945533
$ws.Cells.Item(7, 4).Value = 1.1751505290085449
$ws.Cells.Item(7, 5).Value = 0.58733
$ws.Cells.Item(8, 1).Value = 7
$ws.Cells.Item(8, 2).Value = -271.53604073578464
$ws.Cells.Item(8, 3).Value = 0.04646653013628087
$ws.Cells.Item(8, 4).Value = 4.251964533084473
$ws.Cells.Item(8, 5).Value = 0.0

# ---- Sheet "10": master-problem iteration log ----
$ws = $wb.Worksheets.Item("10")
$ws.Cells.Item(2, 4).Value = 0.04656382242932129
$ws.Cells.Item(2, 5).Value = 70.3188
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = -269.5948197888032
$ws.Cells.Item(3, 3).Value = 0.04513460300575083
$ws.Cells.Item(3, 4).Value = 0.29674409279577635
$ws.Cells.Item(3, 5).Value = 1.35491
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = -268.78676349663795
$ws.Cells.Item(4, 3).Value = 0.09542947205193089
$ws.Cells.Item(4, 4).Value = 0.7397491460584716
$ws.Cells.Item(4, 5).Value = 0.0

